$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New schedule rows for group B2-D2 (rows 284-323), mirrored from the
# existing alternating zebra-stripe style pattern used by rows 2-283.
$newData = @(
    @("Year 5","B2-D2","endocrinology","1","20/12/2025","09:00:00",360),
    @("Year 5","B2-D2","endocrinology","2","21/12/2025","09:00:00",360),
    @("Year 5","B2-D2","endocrinology","3","22/12/2025","09:00:00",360),
    @("Year 5","B2-D2","endocrinology","4","23/12/2025","09:00:00",360),
    @("Year 5","B2-D2","endocrinology","5","24/12/2025","09:00:00",360),
    @("Year 5","B2-D2","endocrinology","6","27/12/2025","09:00:00",360),
    @("Year 5","B2-D2","endocrinology","7","28/12/2025","09:00:00",360),
    @("Year 5","B2-D2","endocrinology","8","29/12/2025","09:00:00",360),
    @("Year 5","B2-D2","endocrinology","9","30/12/2025","09:00:00",360),
    @("Year 5","B2-D2","endocrinology","10","31/12/2025","09:00:00",360),
    @("Year 5","B2-D2","gastroenterology","1","03/01/2026","09:00:00",360),
    @("Year 5","B2-D2","gastroenterology","2","04/01/2026","09:00:00",360),
    @("Year 5","B2-D2","gastroenterology","3","05/01/2026","09:00:00",360),
    @("Year 5","B2-D2","gastroenterology","4","06/01/2026","09:00:00",360),
    @("Year 5","B2-D2","gastroenterology","5","07/01/2026","09:00:00",360),
    @("Year 5","B2-D2","gastroenterology","6","10/01/2026","09:00:00",360),
    @("Year 5","B2-D2","gastroenterology","7","11/01/2026","09:00:00",360),
    @("Year 5","B2-D2","gastroenterology","8","12/01/2026","09:00:00",360),
    @("Year 5","B2-D2","gastroenterology","9","13/01/2026","09:00:00",360),
    @("Year 5","B2-D2","gastroenterology","10","14/01/2026","09:00:00",360),
    @("Year 5","B2-D2","nephrology","1","17/01/2026","09:00:00",360),
    @("Year 5","B2-D2","nephrology","2","18/01/2026","09:00:00",360),
    @("Year 5","B2-D2","nephrology","3","19/01/2026","09:00:00",360),
    @("Year 5","B2-D2","nephrology","4","20/01/2026","09:00:00",360),
    @("Year 5","B2-D2","nephrology","5","21/01/2026","09:00:00",360),
    @("Year 5","B2-D2","neurology","1","06/12/2025","09:00:00",360),
    @("Year 5","B2-D2","neurology","2","07/12/2025","09:00:00",360),
    @("Year 5","B2-D2","neurology","3","08/12/2025","09:00:00",360),
    @("Year 5","B2-D2","neurology","4","09/12/2025","09:00:00",360),
    @("Year 5","B2-D2","neurology","5","13/12/2025","09:00:00",360),
    @("Year 5","B2-D2","neurology","6","14/12/2025","09:00:00",360),
    @("Year 5","B2-D2","neurology","7","15/12/2025","09:00:00",360),
    @("Year 5","B2-D2","neurology","8","16/12/2025","09:00:00",360),
    @("Year 5","B2-D2","physical medicine","1","10/12/2025","09:00:00",360),
    @("Year 5","B2-D2","physical medicine","2","17/12/2025","09:00:00",360),
    @("Year 5","B2-D2","rheumatology","1","07/02/2026","09:00:00",360),
    @("Year 5","B2-D2","rheumatology","2","08/02/2026","09:00:00",360),
    @("Year 5","B2-D2","rheumatology","3","09/02/2026","09:00:00",360),
    @("Year 5","B2-D2","rheumatology","4","10/02/2026","09:00:00",360),
    @("Year 5","B2-D2","rheumatology","5","11/02/2026","09:00:00",360)
)

$startRow = 284
$endRow = $startRow + $newData.Count - 1

# Replicate the alternating row style (fills/number formats for columns
# A-G) from the last existing style pair (rows 282-283, even/odd) down
# across the whole new block in one shot.
$styleSrc = $ws.Range("A282:G283")
$styleDst = $ws.Range("A" + $startRow + ":G" + $endRow)
$styleSrc.Copy($styleDst)

for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $row = $newData[$i]

    if ((($r - 284) % 2) -eq 0) {
        $tmplRow = 282
    } else {
        $tmplRow = 283
    }

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]

    # Columns D ("Session") and E ("Date") hold plain-looking text
    # ("1".."10", "dd/mm/yyyy") that the auto-type-detection would
    # otherwise coerce into a real number / date serial (matching how
    # Excel itself parses typed-in values). The source sheet stores
    # these as literal text, so force text entry via NumberFormat="@"
    # before assigning, then restore the original zebra-stripe style
    # (General / dd-mm-yyyy number format) by pasting formats only from
    # the matching template cell so the visible style is unchanged.
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $row[3]
    $fmtSrc = $ws.Range("D" + $tmplRow)
    $fmtSrc.Copy()
    $ws.Cells.Item($r, 4).PasteSpecial(-4122)

    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 5).Value = $row[4]
    $fmtSrc2 = $ws.Range("E" + $tmplRow)
    $fmtSrc2.Copy()
    $ws.Cells.Item($r, 5).PasteSpecial(-4122)

    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

$excel.CutCopyMode = 0
